$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 121 (shifts existing rows 121:160 down to 122:161)
$ws.Rows("121:121").Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A121").Value = 1
$ws.Range("B121").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C121").Value = "Arica y Parinacota"
$ws.Range("D121").Value = 44489
$ws.Range("E121").Value = 15
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100102
$ws.Range("H121").Value = "Cítricos"
$ws.Range("I121").Value = 100102003
$ws.Range("J121").Value = "Limón"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "2a amarillo"
$ws.Range("M121").Value = 270
$ws.Range("N121").Value = 14000
$ws.Range("O121").Value = 15000
$ws.Range("P121").Value = 14500
$ws.Range("Q121").Value = "$/caja 20 kilos"
$ws.Range("R121").Value = "Región de Coquimbo"
$ws.Range("S121").Value = 725
$ws.Range("T121").Value = 20
